$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.738.80"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "3.800.48"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'596.94"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'167.27"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "'0.160"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").Value = "'35.94"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "4.439.51"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "3.818.90"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "'18.53"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("D17").Value = "67.789.60"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D20").Value = "'461.03"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "'9.90"
$ws.Range("E21").Value = "  -3.28%  "
$ws.Range("D22").Value = "'0.697"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "'83.29"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "'12.06"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'10.01"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("D29").Value = "3.946.30"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").Value = "'7.35"
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'9.04"
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "'3.37"
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "'48.11"
$ws.Range("E43").Value = "  +2.29%  "
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").Value = "'42.80"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").Value = "'147.66"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  +8.14%  "
$ws.Range("D49").Value = "'27.13"
$ws.Range("E49").Value = "  +7.48%  "
$ws.Range("D50").Value = "'394.39"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  +0.91%  "
